$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Time"
$ws.Range("C2").Value = "Student Scheduler"
$ws.Range("D2").Value = "Student Scheduler"
$ws.Range("E2").Value = "Student Scheduler"
$ws.Range("F2").Value = "Student Scheduler"
$ws.Range("G2").Value = "Student Scheduler"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "9:30am - 10:45am"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 3
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "8:00pm - 9:15pm"
$ws.Range("C4").Value = "AVAILABLE"
$ws.Range("D4").Value = "AVAILABLE"
$ws.Range("E4").Value = "AVAILABLE"
$ws.Range("F4").Value = "AVAILABLE"
$ws.Range("G4").Value = "AVAILABLE"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "8:00am - 9:15am"
$ws.Range("C5").Value = "AVAILABLE"
$ws.Range("D5").Value = "AVAILABLE"
$ws.Range("E5").Value = "AVAILABLE"
$ws.Range("F5").Value = "AVAILABLE"
$ws.Range("G5").Value = "AVAILABLE"
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "6:30pm - 7:45pm"
$ws.Range("C6").Value = "AVAILABLE"
$ws.Range("D6").Value = "AVAILABLE"
$ws.Range("E6").Value = "AVAILABLE"
$ws.Range("F6").Value = "AVAILABLE"
$ws.Range("G6").Value = "AVAILABLE"
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "5:00pm - 6:15pm"
$ws.Range("C7").Value = "AVAILABLE"
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = "AVAILABLE"
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = "AVAILABLE"
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "3:30pm - 4:45pm"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = "AVAILABLE"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "2:00pm - 3:15pm"
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = "AVAILABLE"
$ws.Range("E9").Value = "AVAILABLE"
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = "AVAILABLE"
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "12:30pm - 1:45pm"
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = "AVAILABLE"
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 4
$ws.Range("G10").Value = "AVAILABLE"
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "11:00am - 12:15pm"
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = "AVAILABLE"
$ws.Range("E11").Value = 9
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 9
